$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph (and before the "Gameplay" Heading2 paragraph). We do
#    this by collapsing a Range to the very start of the "Gameplay"
#    paragraph and inserting raw OOXML for two paragraphs: the new
#    "Meta description" paragraph followed by the original "Gameplay"
#    heading paragraph (InsertXML replaces the (collapsed) range's
#    owning paragraph, so "Gameplay" must be re-supplied to survive).
# ---------------------------------------------------------------------
$gameplayPara = $d.Paragraphs.Item(2)
$insertRange = $gameplayPara.Range
$insertRange.Collapse(1)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Beat the Beast Griffin’s Gold, a high-volatility slot game with a unique griffin theme. Play for free and experience the enormous winning potential.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Gameplay</w:t></w:r></w:p>'

$insertRange.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Beat the Beast Griffin’s Gold for
#    Free | Review" paragraph that used to sit near the end of the
#    document (right before the italic meta-description paragraph).
#    (Paragraph 1 -- the real H1 title -- has the exact same visible
#    text, so scan every paragraph and keep the *last* match, which is
#    the duplicate near the end of the body.)
# ---------------------------------------------------------------------
$dupTitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd() -eq "Play Beat the Beast Griffin’s Gold for Free | Review") {
        $dupTitlePara = $candidate
    }
}
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the final italic paragraph's text (old meta description)
#    with the new image-generation prompt, preserving its italic run
#    formatting. We assign Range.Text directly (rather than going
#    through Find.Execute's replacement, which smart-quotes straight
#    quotes/apostrophes) so the literal straight quotes/apostrophes in
#    the prompt text come through unchanged.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fullRange = $lastPara.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$newText = "Prompt: Create a feature image for ""Beat the Beast Griffin’s Gold"" that captures the excitement and adventurous spirit of the game. The image should be in a cartoon style and prominently feature a happy Maya warrior with glasses. In the background, include elements such as a temple, columns, and griffin statues to emphasize the game's theme. Add splashes of gold for a touch of luxury and grandeur. The warrior should be depicted holding a crossbow and surrounded by eagles and other mythical creatures to emphasize the game's high volatility and potential for big wins. Make sure the image is visually appealing to attract players and encapsulates the spirit of the game."

$textRange.Text = $newText
